$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Sodium test rows to the new Urine test rows (A2:A5)
$ws.Range("A2").Value = "Urine Color Observation, (CC, SOFTLAB/MIS, COLOR)"
$ws.Range("A3").Value = "Urine Appearance Observation, (CC, SOFTLAB/MIS, APEAR)"
$ws.Range("A4").Value = "Urine Specific Gravity, (CC, SOFTLAB/MIS, SPGR)"
$ws.Range("A5").Value = "pH Measurement Urine Test, (CC, SOFTLAB/MIS, URPH)"

# Fill in the previously-empty rows (A6:A8) with new urine test entries
$ws.Range("A6").Value = "Leukocyte Esterase Urine Test, (CC, SOFTLAB, ULEUK)"
$ws.Range("A7").Value = "Nitrite Urine Test, (CC, SOFTLAB, UNITR)"
$ws.Range("A8").Value = "Protein Qualitative Urine Test, (CC, SOFTLAB/MIS, UPROT)"

# Update the selected cell to A10
$ws.Range("A10").Select()
